$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction in SA algorithm: flatten fitness values for rows 2-250 (C2:C250) to 7293
$ws.Range("C2:C250").Value = 7293
